# "Add files via upload" - refresh the tags list on Sheet1:
#  - shorten the row 47 tag
#  - append the new batch of tag rows the author uploaded (rows 61-90)
#  - leave the selection on the new first empty row below the list

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 47: rename the long tag to the shorter "job titles"
$ws.Range("A47").Value = "job titles"

# New tag rows appended after the existing data (rows 61-90)
$newTags = @(
    "where to get data",
    "job responsibilites",
    "jokes",
    "data science trends",
    "data pipelines",
    "ETL",
    "buying a laptop for data science",
    "data science and software engineering",
    "bye",
    "thanks",
    "comparing to chatgpt",
    "limitations",
    "nosql",
    "books for data science",
    "sql vs ecxel",
    "lablab.me",
    "discord",
    "getting a degree",
    "data modeling",
    "no answer",
    "what are you",
    "options",
    "databases",
    "people in data science",
    "the importance of projects",
    "business intellegnace",
    "AI",
    "BAU",
    "chatbots",
    "chatbots vs LLM"
)

$startRow = 61
for ($i = 0; $i -lt $newTags.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTags[$i]
}

# Move the selection to the first blank row under the refreshed list
$ws.Range("A91").Select()
